$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3 value
$ws.Range("A3").Value = 108999604

# Clear K3 (remove "adult" text) but keep the cell present as an empty string
$ws.Range("K3").Value = ""
